$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for Thbs1-Tnfrsf11b LR pairs (OldD7).
# Target cluster window shifts from [ECs,FAPs,MuSCs] to [FAPs,MuSCs,Resolving-Mac]
# and all dependent statistics are recomputed accordingly.
$updates = @{
    "D2" = "FAPs"
    "G2" = 21.18599966666667
    "H2" = 63.557999
    "I2" = 0.08765141600314529
    "J2" = 0.08765141600314529
    "K2" = 3
    "L2" = 1
    "M2" = 2.214957333333333
    "N2" = 6.644871999999999
    "O2" = 0.8812411509483107
    "P2" = 0.8812411509483107
    "Q2" = 46.92608532568088
    "R2" = 422.334767931128
    "S2" = 0.07724203472086093
    "T2" = 0.07724203472086093
    "D3" = "MuSCs"
    "G3" = 21.18599966666667
    "H3" = 63.557999
    "I3" = 0.08765141600314529
    "J3" = 0.08765141600314529
    "M3" = 0.274148
    "N3" = 0.8224440000000001
    "O3" = 0.1090723037479928
    "P3" = 0.1090723037479928
    "Q3" = 5.808099436617334
    "R3" = 52.272894929556
    "S3" = 0.00956034187023674
    "T3" = 0.00956034187023674
    "D4" = "Resolving-Mac"
    "G4" = 21.18599966666667
    "H4" = 63.557999
    "I4" = 0.08765141600314529
    "J4" = 0.08765141600314529
    "K4" = 1
    "L4" = 0.3333333333333333
    "M4" = 0.02434666666666667
    "N4" = 0.07303999999999999
    "O4" = 0.009686545303696538
    "P4" = 0.009686545303696536
    "Q4" = 0.5158084718844445
    "R4" = 4.64227624696
    "S4" = 0.0008490394120476186
    "T4" = 0.0008490394120476184
    "D5" = "FAPs"
    "I5" = 0.5040014103551328
    "J5" = 0.5040014103551328
    "K5" = 3
    "L5" = 1
    "M5" = 2.214957333333333
    "N5" = 6.644871999999999
    "O5" = 0.8812411509483107
    "P5" = 0.8812411509483107
    "Q5" = 269.8280788269271
    "R5" = 2428.452709442344
    "S5" = 0.4441467829409291
    "T5" = 0.4441467829409291
    "D6" = "MuSCs"
    "I6" = 0.5040014103551328
    "J6" = 0.5040014103551328
    "M6" = 0.274148
    "N6" = 0.8224440000000001
    "O6" = 0.1090723037479928
    "P6" = 0.1090723037479928
    "Q6" = 33.39695399139867
    "R6" = 300.5725859225881
    "S6" = 0.05497259491967182
    "T6" = 0.05497259491967182
    "D7" = "Resolving-Mac"
    "I7" = 0.5040014103551328
    "J7" = 0.5040014103551328
    "K7" = 1
    "L7" = 0.3333333333333333
    "M7" = 0.02434666666666667
    "N7" = 0.07303999999999999
    "O7" = 0.009686545303696538
    "P7" = 0.009686545303696536
    "Q7" = 2.965932658675556
    "R7" = 26.69339392808
    "S7" = 0.004882032494531943
    "T7" = 0.004882032494531942
    "D8" = "FAPs"
    "G8" = 37.20718233333333
    "H8" = 111.621547
    "I8" = 0.1539347809079331
    "J8" = 0.1539347809079331
    "K8" = 3
    "L8" = 1
    "M8" = 2.214957333333333
    "N8" = 6.644871999999999
    "O8" = 0.8812411509483107
    "P8" = 0.8812411509483107
    "Q8" = 82.41232136188709
    "R8" = 741.7108922569839
    "S8" = 0.135653663498283
    "T8" = 0.135653663498283
    "D9" = "MuSCs"
    "G9" = 37.20718233333333
    "H9" = 111.621547
    "I9" = 0.1539347809079331
    "J9" = 0.1539347809079331
    "M9" = 0.274148
    "N9" = 0.8224440000000001
    "O9" = 0.1090723037479928
    "P9" = 0.1090723037479928
    "Q9" = 10.20027462231867
    "R9" = 91.802471600868
    "S9" = 0.01679002118057081
    "T9" = 0.01679002118057081
    "D10" = "Resolving-Mac"
    "G10" = 37.20718233333333
    "H10" = 111.621547
    "I10" = 0.1539347809079331
    "J10" = 0.1539347809079331
    "K10" = 1
    "L10" = 0.3333333333333333
    "M10" = 0.02434666666666667
    "N10" = 0.07303999999999999
    "O10" = 0.009686545303696538
    "P10" = 0.009686545303696536
    "Q10" = 0.9058708658755554
    "R10" = 8.152837792879998
    "S10" = 0.001491096229079295
    "T10" = 0.001491096229079295
    "D11" = "FAPs"
    "G11" = 61.49336899999999
    "H11" = 184.480107
    "I11" = 0.2544123927337887
    "J11" = 0.2544123927337887
    "K11" = 3
    "L11" = 1
    "M11" = 2.214957333333333
    "N11" = 6.644871999999999
    "O11" = 0.8812411509483107
    "P11" = 0.8812411509483107
    "Q11" = 136.2051886179226
    "R11" = 1225.846697561304
    "S11" = 0.2241986697882376
    "T11" = 0.2241986697882376
    "D12" = "MuSCs"
    "G12" = 61.49336899999999
    "H12" = 184.480107
    "I12" = 0.2544123927337887
    "J12" = 0.2544123927337887
    "M12" = 0.274148
    "N12" = 0.8224440000000001
    "O12" = 0.1090723037479928
    "P12" = 0.1090723037479928
    "Q12" = 16.858284124612
    "R12" = 151.724557121508
    "S12" = 0.02774934577751344
    "T12" = 0.02774934577751344
    "D13" = "Resolving-Mac"
    "G13" = 61.49336899999999
    "H13" = 184.480107
    "I13" = 0.2544123927337887
    "J13" = 0.2544123927337887
    "K13" = 1
    "L13" = 0.3333333333333333
    "M13" = 0.02434666666666667
    "N13" = 0.07303999999999999
    "O13" = 0.009686545303696538
    "P13" = 0.009686545303696536
    "Q13" = 1.497158557253333
    "R13" = 13.47442701528
    "S13" = 0.00246437716803768
    "T13" = 0.00246437716803768
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
